$wb = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")
$ws = $wb.Worksheets.Item("BGDPbES")

# "BAU Guaranteed Dispatch" input (column B) switches from 0 to 1 for the
# following electricity sources (rows on the BGDPbES sheet). The rest of each
# row (columns C:AK) already holds formulas that reference column B, so they
# recalc automatically once B changes.
#   Row 6  = onshore wind
#   Row 7  = solar PV
#   Row 8  = solar thermal
#   Row 9  = biomass
#   Row 10 = geothermal
#   Row 14 = offshore wind
#   Row 17 = municipal solid waste
$rows = @(6, 7, 8, 9, 10, 14, 17)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 2).Value = 1
}

# Row 17 (municipal solid waste) previously held hardcoded 0 values across
# C17:AK17 instead of formulas referencing column B like every other row in
# the table. Bring it in line with the rest of the sheet by making C17:AK17
# reference $B$17.
$ws.Range("C17").Formula = "=`$B`$17"
$ws.Range("D17:AK17").Formula = "=`$B`$17"

# Recalculate so the cached formula results reflect the new inputs.
$excel.CalculateFull()

# Restore the on-screen selections left behind after making the edits:
# the BGDPbES sheet was scrolled/selected around the edited row 17 range,
# then focus returned to the About sheet before saving.
$ws.Activate()
$ws.Range("C17:AK17").Select()

$wsAbout.Activate()
$wsAbout.Range("C14").Select()
